$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting rows 207:411 down to 208:412
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across the whole dataset.
$ws.Range("A207").Value = 3
$ws.Range("B207").Value = "Femacal de La Calera"
$ws.Range("C207").Value = "Coquimbo"
$ws.Range("D207").Value = 45271
$ws.Range("E207").Value = 5
$ws.Range("F207").Value = 100112039
$ws.Range("G207").Value = "Ciboulette"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 120
$ws.Range("K207").Value = 1500
$ws.Range("L207").Value = 1500
$ws.Range("M207").Value = 1500
$ws.Range("N207").Value = "`$/docena de atados"
$ws.Range("O207").Value = "Provincia de Quillota"
$ws.Range("P207").Value = 500
$ws.Range("Q207").Value = 3
$ws.Range("R207").Value = "Hortaliza"
